$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 39700
$ws.Range("D2").Value = 57361868
$ws.Range("C3").Value = 94928
$ws.Range("D3").Value = 139109833
$ws.Range("C4").Value = 32374
$ws.Range("D4").Value = 47929001
$ws.Range("C5").Value = 9136
$ws.Range("D5").Value = 13577411
$ws.Range("C6").Value = 2163
$ws.Range("D6").Value = 3212648
$ws.Range("C7").Value = 194
$ws.Range("D7").Value = 286093
$ws.Range("C12").Value = 43067
$ws.Range("D12").Value = 58360363
$ws.Range("C13").Value = 10108
$ws.Range("D13").Value = 14605861
$ws.Range("C14").Value = 26934
$ws.Range("D14").Value = 39479031
$ws.Range("C15").Value = 8587
$ws.Range("D15").Value = 12743534
$ws.Range("C16").Value = 2246
$ws.Range("D16").Value = 3336653
$ws.Range("C17").Value = 444
$ws.Range("D17").Value = 655123
$ws.Range("C20").Value = 10600
$ws.Range("D20").Value = 13997007
$ws.Range("C21").Value = 13952
$ws.Range("D21").Value = 20125542
$ws.Range("C22").Value = 32783
$ws.Range("D22").Value = 48083857
$ws.Range("C23").Value = 10561
$ws.Range("D23").Value = 15695505
$ws.Range("C24").Value = 2751
$ws.Range("D24").Value = 4091174
$ws.Range("C25").Value = 562
$ws.Range("D25").Value = 837092
$ws.Range("C26").Value = 41
$ws.Range("D26").Value = 60953
$ws.Range("C27").Value = 12119
$ws.Range("D27").Value = 16152288
$ws.Range("C28").Value = 8067
$ws.Range("D28").Value = 11668113
$ws.Range("C29").Value = 23417
$ws.Range("D29").Value = 34365716
$ws.Range("C30").Value = 8064
$ws.Range("D30").Value = 11988633
$ws.Range("C31").Value = 2044
$ws.Range("D31").Value = 3049199
$ws.Range("C32").Value = 388
$ws.Range("D32").Value = 579415
$ws.Range("C33").Value = 32
$ws.Range("D33").Value = 47893
$ws.Range("C34").Value = 8656
$ws.Range("D34").Value = 11430299
$ws.Range("C35").Value = 3440
$ws.Range("D35").Value = 4966498
$ws.Range("C36").Value = 8220
$ws.Range("D36").Value = 12008499
$ws.Range("C37").Value = 3285
$ws.Range("D37").Value = 4870461
$ws.Range("C38").Value = 851
$ws.Range("D38").Value = 1267555
$ws.Range("C39").Value = 173
$ws.Range("D39").Value = 257186
$ws.Range("C40").Value = 7
$ws.Range("D40").Value = 10500
$ws.Range("C41").Value = 2601
$ws.Range("D41").Value = 3514777
$ws.Range("C42").Value = 18059
$ws.Range("D42").Value = 26087869
$ws.Range("C43").Value = 53075
$ws.Range("D43").Value = 77771397
$ws.Range("C44").Value = 19538
$ws.Range("D44").Value = 29009453
$ws.Range("C45").Value = 5847
$ws.Range("D45").Value = 8702435
$ws.Range("C46").Value = 1303
$ws.Range("D46").Value = 1944144
$ws.Range("C47").Value = 74
$ws.Range("D47").Value = 109015
$ws.Range("C50").Value = 17472
$ws.Range("D50").Value = 23178144
$ws.Range("C51").Value = 2221
$ws.Range("D51").Value = 3224544
$ws.Range("C52").Value = 7481
$ws.Range("D52").Value = 10992497
$ws.Range("C53").Value = 2491
$ws.Range("D53").Value = 3718184
$ws.Range("C54").Value = 785
$ws.Range("D54").Value = 1172415
$ws.Range("C57").Value = 7570
$ws.Range("D57").Value = 10413719
$ws.Range("C58").Value = 1330
$ws.Range("D58").Value = 2421499
$ws.Range("C59").Value = 3287
$ws.Range("D59").Value = 6014557
$ws.Range("C60").Value = 1288
$ws.Range("D60").Value = 2361066
$ws.Range("C61").Value = 437
$ws.Range("D61").Value = 798083
$ws.Range("C62").Value = 149
$ws.Range("D62").Value = 284600
$ws.Range("C64").Value = 2017
$ws.Range("D64").Value = 3402026
$ws.Range("C65").Value = 16195
$ws.Range("D65").Value = 23385183
$ws.Range("C66").Value = 46617
$ws.Range("D66").Value = 68172221
$ws.Range("C67").Value = 16269
$ws.Range("D67").Value = 24171049
$ws.Range("C68").Value = 4751
$ws.Range("D68").Value = 7076288
$ws.Range("C69").Value = 1006
$ws.Range("D69").Value = 1495329
$ws.Range("C70").Value = 84
$ws.Range("D70").Value = 123330
$ws.Range("C73").Value = 15657
$ws.Range("D73").Value = 20592228
$ws.Range("C74").Value = 57149
$ws.Range("D74").Value = 83108354
$ws.Range("C75").Value = 158216
$ws.Range("D75").Value = 232962636
$ws.Range("C76").Value = 67896
$ws.Range("D76").Value = 101139835
$ws.Range("C77").Value = 21845
$ws.Range("D77").Value = 32640013
$ws.Range("C78").Value = 5298
$ws.Range("D78").Value = 7913102
$ws.Range("C79").Value = 323
$ws.Range("D79").Value = 479670
$ws.Range("C82").Value = 4
$ws.Range("D82").Value = 6000
$ws.Range("C83").Value = 5
$ws.Range("D83").Value = 7500
$ws.Range("C85").Value = 56041
$ws.Range("D85").Value = 75935633
$ws.Range("C86").Value = 4900
$ws.Range("D86").Value = 7101061
$ws.Range("C87").Value = 12149
$ws.Range("D87").Value = 17845074
$ws.Range("C88").Value = 4022
$ws.Range("D88").Value = 5992958
$ws.Range("C89").Value = 1395
$ws.Range("D89").Value = 2083611
$ws.Range("C90").Value = 308
$ws.Range("D90").Value = 459512
$ws.Range("C91").Value = 30
$ws.Range("D91").Value = 44902
$ws.Range("C93").Value = 5687
$ws.Range("D93").Value = 7640738
$ws.Range("C94").Value = 1712
$ws.Range("D94").Value = 2467363
$ws.Range("C95").Value = 5515
$ws.Range("D95").Value = 8126855
$ws.Range("C96").Value = 2025
$ws.Range("D96").Value = 3014308
$ws.Range("C97").Value = 738
$ws.Range("D97").Value = 1105960
$ws.Range("C101").Value = 3775
$ws.Range("D101").Value = 5003552
$ws.Range("C102").Value = 791
$ws.Range("D102").Value = 1408715
$ws.Range("C103").Value = 495
$ws.Range("D103").Value = 914604
$ws.Range("C104").Value = 191
$ws.Range("D104").Value = 356289
$ws.Range("C105").Value = 63
$ws.Range("D105").Value = 115500
$ws.Range("C107").Value = 11354
$ws.Range("D107").Value = 16463626
$ws.Range("C108").Value = 30234
$ws.Range("D108").Value = 44395520
$ws.Range("C109").Value = 10129
$ws.Range("D109").Value = 15059226
$ws.Range("C110").Value = 2789
$ws.Range("D110").Value = 4158383
$ws.Range("C111").Value = 523
$ws.Range("D111").Value = 778903
$ws.Range("C114").Value = 10169
$ws.Range("D114").Value = 13403430
$ws.Range("C115").Value = 31934
$ws.Range("D115").Value = 46026852
$ws.Range("C116").Value = 68648
$ws.Range("D116").Value = 100432607
$ws.Range("C117").Value = 22054
$ws.Range("D117").Value = 32767662
$ws.Range("C118").Value = 6286
$ws.Range("D118").Value = 9359011
$ws.Range("C119").Value = 1196
$ws.Range("D119").Value = 1786889
$ws.Range("C120").Value = 99
$ws.Range("D120").Value = 144895
$ws.Range("C124").Value = 26782
$ws.Range("D124").Value = 35722593
$ws.Range("C125").Value = 37920
$ws.Range("D125").Value = 54694577
$ws.Range("C126").Value = 80129
$ws.Range("D126").Value = 117135640
$ws.Range("C127").Value = 24729
$ws.Range("D127").Value = 36701297
$ws.Range("C128").Value = 6644
$ws.Range("D128").Value = 9873243
$ws.Range("C129").Value = 1341
$ws.Range("D129").Value = 1989958
$ws.Range("C130").Value = 73
$ws.Range("D130").Value = 107728
$ws.Range("C133").Value = 33094
$ws.Range("D133").Value = 43891256
$ws.Range("C134").Value = 13941
$ws.Range("D134").Value = 20175447
$ws.Range("C135").Value = 33575
$ws.Range("D135").Value = 49294238
$ws.Range("C136").Value = 11850
$ws.Range("D136").Value = 17607946
$ws.Range("C137").Value = 3110
$ws.Range("D137").Value = 4636241
$ws.Range("C138").Value = 532
$ws.Range("D138").Value = 791990
$ws.Range("C141").Value = 11231
$ws.Range("D141").Value = 14955693
$ws.Range("C142").Value = 37045
$ws.Range("D142").Value = 53497150
$ws.Range("C143").Value = 85240
$ws.Range("D143").Value = 124859507
$ws.Range("C144").Value = 25361
$ws.Range("D144").Value = 37677796
$ws.Range("C145").Value = 6670
$ws.Range("D145").Value = 9950937
$ws.Range("C146").Value = 1525
$ws.Range("D146").Value = 2266802
$ws.Range("C147").Value = 92
$ws.Range("D147").Value = 137630
$ws.Range("C149").Value = 30479
$ws.Range("D149").Value = 41052060
